$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column N (2020) data
$ws.Range("N4").Value = 2020
$ws.Range("N5").Value = 534

# Update sheet view: scroll and selection
$ws.Application.ActiveWindow.ScrollColumn = 5
$ws.Range("S10").Select()
